$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.935.86"
$ws.Range("E2").Value = "  -1.68%  "
$ws.Range("D3").Value = "1.834.83"
$ws.Range("E3").Value = "  -1.89%  "
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'245.43"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").Value = "'0.6911"
$ws.Range("E6").Value = "  -2.08%  "
$ws.Range("D7").Value = "'0.9994"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'0.07694"
$ws.Range("E8").Value = "  -2.68%  "
$ws.Range("D9").Value = "'0.3051"
$ws.Range("E9").Value = "  -2.83%  "
$ws.Range("D10").Value = "'23.59"
$ws.Range("E10").Value = "  -4.07%  "
$ws.Range("D11").Value = "'0.07805"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("D12").Value = "1.838.90"
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("D13").Value = "'5.080"
$ws.Range("E13").Value = "  -2.39%  "
$ws.Range("D14").Value = "'90.64"
$ws.Range("E14").Value = "  -3.52%  "
$ws.Range("D15").Value = "'0.6815"
$ws.Range("E15").Value = "  -3.03%  "
$ws.Range("D16").Value = "'6.448"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").Value = "'0.000008366"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "28.921.31"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("D19").Value = "'243.46"
$ws.Range("E19").Value = "  -4.35%  "
$ws.Range("D20").Value = "2.082.14"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("D21").Value = "'12.70"
$ws.Range("E21").Value = "  -3.29%  "
$ws.Range("D22").Value = "'0.9992"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").Value = "'7.484"
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").Value = "'0.1470"
$ws.Range("E25").Value = "  -5.87%  "
$ws.Range("D26").Value = "'161.90"
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").Value = "'8.817"
$ws.Range("E27").Value = "  -2.22%  "
$ws.Range("D28").Value = "'18.21"
$ws.Range("E28").Value = "  -3.22%  "
$ws.Range("D29").Value = "'1.554"
$ws.Range("E29").Value = "  +3.06%  "
$ws.Range("D30").Value = "'4.218"
$ws.Range("E30").Value = "  -2.92%  "
$ws.Range("D31").Value = "'4.161"
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("D32").Value = "'1.179"
$ws.Range("E32").Value = "  -2.41%  "
$ws.Range("D33").Value = "'0.05133"
$ws.Range("E33").Value = "  -3.30%  "
$ws.Range("D34").Value = "'0.7682"
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("D35").Value = "'1.850"
$ws.Range("E35").Value = "  -2.64%  "
$ws.Range("D36").Value = "'1.148"
$ws.Range("E36").Value = "  -2.48%  "
$ws.Range("D37").Value = "'2.682"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("E38").Value = "  -2.48%  "
$ws.Range("D39").Value = "1.231.15"
$ws.Range("E39").Value = "  -4.01%  "
$ws.Range("D40").Value = "'2.694"
$ws.Range("E40").Value = "  -2.65%  "
$ws.Range("D41").Value = "'0.9219"
$ws.Range("E41").Value = "  +2.94%  "
$ws.Range("D42").Value = "'108.36"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("D43").Value = "'5.851"
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("D44").Value = "'0.9990"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "'9.634"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").Value = "1.982.21"
$ws.Range("E46").Value = "  -2.48%  "
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("E48").Value = "  -4.64%  "
$ws.Range("D49").Value = "'64.31"
$ws.Range("E49").Value = "  -9.95%  "
$ws.Range("D50").Value = "'1.749"
$ws.Range("E50").Value = "  -2.82%  "
$ws.Range("D51").Value = "'6.935"
$ws.Range("E51").Value = "  -2.03%  "
